# dir_Area.xlsx update ("Add files via upload")
#
# - Header row (row 1) relabelled: A1 "AREA" -> "Valor", B1 "Ciudad" -> "Categoría",
#   and both header cells get centered horizontal alignment.
# - Row 25 city name typo fixed: "San André" -> "San Andrés".
# - Selection reset to A1 (was sitting on C28 in the source file).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo first so new shared-string entries land in the same order
# the workbook was edited in (San Andrés, then the new headers).
$ws.Range("B25").Value = "San Andrés"

# Re-label the header row.
$ws.Range("A1").Value = "Valor"
$ws.Range("B1").Value = "Categoría"

# Center the header cells (adds a new cell style picked up by A1:B1).
$ws.Range("A1:B1").HorizontalAlignment = -4108

# Clear the stray selection that was left on C28 -> put it back on A1.
[void]$ws.Range("A1").Select()
